$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Test Date" header in B1
$ws.Range("B1").Value = "Test Date"

# Add date value in B2 (serial 45588 == 2024-10-23), formatted as a (built-in) date
$ws.Range("B2").Value = 45588
$ws.Range("B2").NumberFormat = "mm-dd-yy"

# Size column B to fit the date values (stored width ends up 10.5 in the XML)
$ws.Columns.Item(2).ColumnWidth = 9.6666666666666667

# Update selection to B3, matching where the user left the cursor after entry
$ws.Range("B3").Select()
